# loss made pretty 2
#
# Before:
#   B1 = 0 (numeric, styled)
#   A2 = "Training" (styled), B2 = [training losses]
#   A3 = "Validation" (styled), B3 = [validation losses]
#
# After:
#   B1 = "Training" (styled), C1 = "Validation" (styled)
#   A2 = 0 (numeric, styled), B2 = [new training losses], C2 = [new validation losses]

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Propagate the existing header style (currently on B1) onto the cells
#    that will need it in the new layout (A2 and C1) before we touch any
#    values, so the style index is reused rather than a new one created.
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# 2) The old row 3 ("Validation" label + its values) is no longer needed;
#    removing it also shrinks the used range back down to 2 rows.
$ws.Rows("3:3").Delete()

# 3) Write the new cell contents.
$ws.Range("A2").Value = 0
$ws.Range("B1").Value = "Training"
$ws.Range("C1").Value = "Validation"
$ws.Range("B2").Value = "[0.03761576488614082, 0.007720354124903679, 0.004869405888020992, 0.004016904495656491, 0.003860360607504845, 0.0033066946268081663, 0.003209751956164837, 0.003131478950381279, 0.002889433428645134, 0.002608037553727627]"
$ws.Range("C2").Value = "[4.466980361938477, 0.013171311095356941, 0.013815078884363174, 0.018883594870567323, 0.024239375442266464, 0.022170404344797133, 0.02061031311750412, 0.02529122605919838, 0.02824832797050476, 0.02674718052148819]"
